# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Fri Jun 28 11:19:43 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.543.30"
$ws.Range("E2").Value = "'  +0.60%  "
$ws.Range("D3").Value = "'3.445.67"
$ws.Range("E3").Value = "'  +1.33%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'577.01"
$ws.Range("E5").Value = "'  +0.58%  "
$ws.Range("D6").Value = "'145.30"
$ws.Range("E6").Value = "'  +4.47%  "
$ws.Range("D7").Value = "'3.449.12"
$ws.Range("E7").Value = "'  +1.42%  "
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("E9").Value = "'  +1.82%  "
$ws.Range("E10").Value = "'  -0.23%  "
$ws.Range("E11").Value = "'  +3.50%  "
$ws.Range("D12").Value = "'0.391"
$ws.Range("E12").Value = "'  +2.58%  "
$ws.Range("D13").Value = "'4.032.49"
$ws.Range("E13").Value = "'  +1.35%  "
$ws.Range("D14").Value = "'28.42"
$ws.Range("E14").Value = "'  +6.86%  "
$ws.Range("E15").Value = "'  -0.36%  "
$ws.Range("E16").Value = "'  +0.91%  "
$ws.Range("D17").Value = "'3.442.35"
$ws.Range("E17").Value = "'  +1.34%  "
$ws.Range("D18").Value = "'61.624.90"
$ws.Range("E18").Value = "'  +0.70%  "
$ws.Range("D19").Value = "'6.38"
$ws.Range("E19").Value = "'  +7.45%  "
$ws.Range("D20").Value = "'14.31"
$ws.Range("E20").Value = "'  +3.28%  "
$ws.Range("E21").Value = "'  +1.03%  "
$ws.Range("D22").Value = "'401.81"
$ws.Range("E22").Value = "'  +7.00%  "
$ws.Range("D23").Value = "'0.568"
$ws.Range("E23").Value = "'  +2.77%  "
$ws.Range("D24").Value = "'74.45"
$ws.Range("E24").Value = "'  +4.69%  "
$ws.Range("E25").Value = "'  +0.56%  "
$ws.Range("E26").Value = "'  -0.47%  "
$ws.Range("D27").Value = "'0.0000123"
$ws.Range("E27").Value = "'  +0.87%  "
$ws.Range("D28").Value = "'3.588.55"
$ws.Range("E28").Value = "'  +1.59%  "
$ws.Range("E29").Value = "'  +4.28%  "
$ws.Range("D30").Value = "'7.61"
$ws.Range("E30").Value = "'  +2.55%  "
$ws.Range("E31").Value = "'  -0.05%  "
$ws.Range("D32").Value = "'8.25"
$ws.Range("E32").Value = "'  +1.44%  "
$ws.Range("E33").Value = "'  +1.95%  "
$ws.Range("E34").Value = "'  -10.53%  "
$ws.Range("E35").Value = "'  -0.07%  "
$ws.Range("D36").Value = "'23.94"
$ws.Range("E36").Value = "'  +1.99%  "
$ws.Range("E37").Value = "'  +2.40%  "
$ws.Range("D38").Value = "'3.471.48"
$ws.Range("E38").Value = "'  +1.57%  "
$ws.Range("B39").Value = "'ImmutableX"
$ws.Range("C39").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "'  -0.35%  "
$ws.Range("B40").Value = "'NEARProtocol"
$ws.Range("C40").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'5.13"
$ws.Range("E40").Value = "'  +0.21%  "
$ws.Range("D41").Value = "'167.11"
$ws.Range("E41").Value = "'  +0.46%  "
$ws.Range("D42").Value = "'0.0791"
$ws.Range("E42").Value = "'  +2.56%  "
$ws.Range("E43").Value = "'  +5.00%  "
$ws.Range("E44").Value = "'  +3.13%  "
$ws.Range("D45").Value = "'4.52"
$ws.Range("E45").Value = "'  +2.69%  "
$ws.Range("E46").Value = "'  -0.27%  "
$ws.Range("B47").Value = "'FirstDigitalUSD"
$ws.Range("C47").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "'  -0.02%  "
$ws.Range("B48").Value = "'OKB"
$ws.Range("C48").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'42.42"
$ws.Range("E48").Value = "'  +1.13%  "
$ws.Range("D49").Value = "'2.611.58"
$ws.Range("E49").Value = "'  +3.74%  "
$ws.Range("E50").Value = "'  -2.02%  "
$ws.Range("E51").Value = "'  +2.49%  "

Write-Output "Applied 87 cell updates"
